# WebForm User Assignment execution
# Update the PN_Value (phone number) column F for rows 2-10 with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = "9840012393"
    3  = "9840047034"
    4  = "9840011675"
    5  = "9840014446"
    6  = "9840083166"
    7  = "9840023824"
    8  = "9840004860"
    9  = "9840059068"
    10 = "9840036323"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 6)
    # Force the number-looking phone number to be stored as text (matching
    # the original column's data type) instead of being auto-converted to a
    # numeric value, then restore the original (default/"Normal") cell style.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
    $cell.Style = "Normal"
}
